$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert "WEO" / "Weapon Emitter Overdrive" row, alphabetically between
# "WAQT" (row 265) and "WSS" (row 266).
$ws.Rows.Item(266).EntireRow.Insert()
$ws.Range("A266").Value = "WEO"
$ws.Range("B266").Value = "Weapon Emitter Overdrive"
$ws.Range("C266").Value = "Starship Trait (Infinity Promotional Ship)"

# Insert "TG" / "Terran Goodbye" row, alphabetically between "TFO" (row 247)
# and "ThS" (row 248).
$ws.Rows.Item(248).EntireRow.Insert()
$ws.Range("A248").Value = "TG"
$ws.Range("B248").Value = "Terran Goodbye"
$ws.Range("C248").Value = "Starship Trait (Lockbox Ship)"

$ws.Range("A248:XFD248").Select()
